$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 6.162999999999999
$ws.Range("B18").Value = 5.137
$ws.Range("B20").Value = 6.862
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("B35").Value = 8.300999999999998
$ws.Range("B69").Value = 5.627
$ws.Range("B76").Value = 6.544999999999999
$ws.Range("B78").Value = 8.221
$ws.Range("B82").Value = 5.457
$ws.Range("B83").Value = 5.88
$ws.Range("B93").Value = 4.973000000000001
